$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Input data changes (row 5: Månedsinntekt per arbeidsgiver) ---
$ws.Range("C5").Value = 31000
$ws.Range("D5").Value = 31000
$ws.Range("E5").ClearContents()

# --- Row 13 (Ønsket refusjon): turn formulas into plain input values ---
$ws.Range("C13").Value = 15000
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0

# --- Row 28 (Dekningsgrad) ---
$ws.Range("C28").Value = 0.8
$ws.Range("D28").Value = 0.6
$ws.Range("E28").Value = 1
$ws.Range("F28").NumberFormat = $ws.Range("E28").NumberFormat
$ws.Range("F29").NumberFormat = $ws.Range("E29").NumberFormat

# --- Row 24: formula changed from "X10*X17/$B$22" to "12*X16/$B$14" ---
$ws.Range("C24").Formula = '=12*C16/$B$14'
$ws.Range("D24").Formula = '=12*D16/$B$14'
$ws.Range("E24").Formula = '=12*E16/$B$14'

# --- Row 54: fix self-referencing typo bug (C52 -> C54, $C$52:$D$52 -> $C$54:$D$54) ---
$ws.Range("D54").Formula = '=IF(($B$53-C54)>0,IF(LARGE($C$52:$F$52,$B$47)<=D$52,1,0), 0)'
$ws.Range("E54").Formula = '=IF(($B$53-SUM($C$54:$D$54))>0,IF(LARGE($C$52:$F$52,$B$47)<=E$52,1,0), 0)'

# --- sheet view: selection moved from B36 to B38 (this also drops topLeftCell) ---
$ws.Range("B38").Select()

$wb.Application.Calculate()
